$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'64.613.66"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.425.53"
$ws.Range("E3").Value = "  -1.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'573.75"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'157.12"
$ws.Range("E6").Value = "  -2.26%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  +7.02%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'3.426.62"
$ws.Range("E9").Value = "  -1.09%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -2.85%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.124"
$ws.Range("E11").Value = "  -1.57%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.56%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'4.011.76"
$ws.Range("E13").Value = "  -1.35%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.29%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  -3.36%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'27.92"
$ws.Range("E16").Value = "  -0.81%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'64.610.94"
$ws.Range("E17").Value = "  -0.63%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.405.58"
$ws.Range("E18").Value = "  -1.77%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "  +0.18%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'13.98"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'378.65"
$ws.Range("E21").Value = "  -2.61%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.33%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.18%  "

# Row 24 - now Litecoin (was Dai)
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'72.51"
$ws.Range("E24").Value = "  -0.90%  "

# Row 25 - now Dai (was Litecoin)
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.22%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  -3.13%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'10.29"
$ws.Range("E27").Value = "  +5.96%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -1.25%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - Fetch.AI
$ws.Range("D30").Value = "'1.49"
$ws.Range("E30").Value = "  +3.41%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'6.23"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -0.55%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'23.17"
$ws.Range("E33").Value = "  -2.12%  "

# Row 34 - Aptos
$ws.Range("D34").Value = "'7.25"
$ws.Range("E34").Value = "  +2.33%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +7.34%  "

# Row 36 - Monero
$ws.Range("D36").Value = "'159.59"
$ws.Range("E36").Value = "  -2.24%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -1.27%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "'6.98"
$ws.Range("E38").Value = "  +6.67%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.0764"
$ws.Range("E39").Value = "  -0.41%  "

# Row 40 - EnergySwap
$ws.Range("D40").Value = "'26.88"
$ws.Range("E40").Value = "  -1.28%  "

# Row 41 - Maker
$ws.Range("D41").Value = "'2.878.79"
$ws.Range("E41").Value = "  -4.06%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +1.49%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "'26.69"
$ws.Range("E43").Value = "  +9.30%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "'0.0318"
$ws.Range("E44").Value = "  +0.58%  "

# Row 45 - OKB
$ws.Range("D45").Value = "'42.84"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "'0.773"
$ws.Range("E46").Value = "  -0.31%  "

# Row 47 - Bittensor
$ws.Range("D47").Value = "'320.08"
$ws.Range("E47").Value = "  +4.78%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "'1.09"
$ws.Range("E48").Value = "  -0.52%  "

# Row 49 - Stellar
$ws.Range("D49").Value = "'0.109"
$ws.Range("E49").Value = "  +1.92%  "

# Row 50 - dogwifhat
$ws.Range("E50").Value = "  +1.02%  "

# Row 51 - Cosmos
$ws.Range("D51").Value = "'6.58"
$ws.Range("E51").Value = "  +0.10%  "
